$wb = $excel.ActiveWorkbook

# This script applies a scheduled market-data refresh to the per-job profit
# sheets (currentAveragePrice / LevePrice / LeveProfit columns, H:N) across
# several worksheets. All values are static data (no formulas in this workbook).

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3523.4285
$ws.Range("J64").Value = 3888.889
$ws.Range("L64").Value = 3888.889
$ws.Range("N64").Value = -4384.889
$ws.Range("H67").Value = 3523.4285
$ws.Range("J67").Value = 3888.889
$ws.Range("L67").Value = 3888.889
$ws.Range("N67").Value = -5604.889
$ws.Range("H80").Value = 395.7143
$ws.Range("I80").Value = 479
$ws.Range("J80").Value = 362.4
$ws.Range("K80").Value = 1437
$ws.Range("L80").Value = 1087.2
$ws.Range("M80").Value = -439
$ws.Range("N80").Value = -3083.2
$ws.Range("H82").Value = 141.16667
$ws.Range("I82").Value = 141.16667
$ws.Range("K82").Value = 423.50001
$ws.Range("M82").Value = -17.50001000000003
$ws.Range("H83").Value = 395.7143
$ws.Range("I83").Value = 479
$ws.Range("J83").Value = 362.4
$ws.Range("K83").Value = 4311
$ws.Range("L83").Value = 3261.6
$ws.Range("M83").Value = 681
$ws.Range("N83").Value = -13245.6
$ws.Range("H85").Value = 141.16667
$ws.Range("I85").Value = 141.16667
$ws.Range("K85").Value = 423.50001
$ws.Range("M85").Value = 980.49999
$ws.Range("H92").Value = 1113.0769
$ws.Range("I92").Value = 147
$ws.Range("J92").Value = 4333.3335
$ws.Range("K92").Value = 147
$ws.Range("L92").Value = 4333.3335
$ws.Range("M92").Value = 1101
$ws.Range("N92").Value = -6829.3335
$ws.Range("H95").Value = 18980
$ws.Range("J95").Value = 18980
$ws.Range("L95").Value = 18980
$ws.Range("N95").Value = -24472
$ws.Range("H138").Value = 2634.0527
$ws.Range("I138").Value = 2110.4443
$ws.Range("J138").Value = 2796.5518
$ws.Range("K138").Value = 6331.3329
$ws.Range("L138").Value = 8389.6554
$ws.Range("M138").Value = -1191.3329
$ws.Range("N138").Value = -18669.6554

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11909997
$ws.Range("J32").Value = 20662
$ws.Range("L32").Value = 20662
$ws.Range("N32").Value = -21236

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2087.6843
$ws.Range("I86").Value = 1705.0769
$ws.Range("J86").Value = 2916.6667
$ws.Range("K86").Value = 1705.0769
$ws.Range("L86").Value = 2916.6667
$ws.Range("M86").Value = -582.0769
$ws.Range("N86").Value = -5162.6667
$ws.Range("H89").Value = 2087.6843
$ws.Range("I89").Value = 1705.0769
$ws.Range("J89").Value = 2916.6667
$ws.Range("K89").Value = 8525.3845
$ws.Range("L89").Value = 14583.3335
$ws.Range("M89").Value = -2909.3845
$ws.Range("N89").Value = -25815.3335

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H103").Value = 14592.4
$ws.Range("I103").Value = 10987.333
$ws.Range("K103").Value = 10987.333
$ws.Range("M103").Value = -9815.333000000001

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 563.6
$ws.Range("I5").Value = 559.3570999999999
$ws.Range("J5").Value = 567.3125
$ws.Range("K5").Value = 1678.0713
$ws.Range("L5").Value = 1701.9375
$ws.Range("M5").Value = -1566.0713
$ws.Range("N5").Value = -1925.9375
$ws.Range("H57").Value = 3211.111
$ws.Range("I57").Value = 966.6667
$ws.Range("J57").Value = 4333.3335
$ws.Range("K57").Value = 2900.0001
$ws.Range("L57").Value = 13000.0005
$ws.Range("M57").Value = -2341.0001
$ws.Range("N57").Value = -14118.0005
$ws.Range("H58").Value = 7800
$ws.Range("J58").Value = 7800
$ws.Range("L58").Value = 23400
$ws.Range("N58").Value = -23656
$ws.Range("H59").Value = 1091.5
$ws.Range("I59").Value = 283
$ws.Range("K59").Value = 849
$ws.Range("M59").Value = -309
$ws.Range("H60").Value = 320.69232
$ws.Range("I60").Value = 226.25
$ws.Range("J60").Value = 362.66666
$ws.Range("K60").Value = 678.75
$ws.Range("L60").Value = 1087.99998
$ws.Range("M60").Value = -427.75
$ws.Range("N60").Value = -1589.99998
$ws.Range("H68").Value = 1136.129
$ws.Range("I68").Value = 1075.2
$ws.Range("J68").Value = 1246.909
$ws.Range("K68").Value = 3225.6
$ws.Range("L68").Value = 3740.727
$ws.Range("M68").Value = -2414.6
$ws.Range("N68").Value = -5362.727000000001
$ws.Range("H69").Value = 3055.1428
$ws.Range("I69").Value = 900
$ws.Range("J69").Value = 3220.923
$ws.Range("K69").Value = 2700
$ws.Range("L69").Value = 9662.769
$ws.Range("M69").Value = -1889
$ws.Range("N69").Value = -11284.769
$ws.Range("H70").Value = 3550
$ws.Range("I70").Value = 1500
$ws.Range("J70").Value = 3960
$ws.Range("K70").Value = 4500
$ws.Range("L70").Value = 11880
$ws.Range("M70").Value = -4185
$ws.Range("N70").Value = -12510
$ws.Range("H71").Value = 1136.129
$ws.Range("I71").Value = 1075.2
$ws.Range("J71").Value = 1246.909
$ws.Range("K71").Value = 9676.800000000001
$ws.Range("L71").Value = 11222.181
$ws.Range("M71").Value = -5620.800000000001
$ws.Range("N71").Value = -19334.181
$ws.Range("H72").Value = 3055.1428
$ws.Range("I72").Value = 900
$ws.Range("J72").Value = 3220.923
$ws.Range("K72").Value = 8100
$ws.Range("L72").Value = 28988.307
$ws.Range("M72").Value = -4044
$ws.Range("N72").Value = -37100.307
$ws.Range("H73").Value = 3550
$ws.Range("I73").Value = 1500
$ws.Range("J73").Value = 3960
$ws.Range("K73").Value = 4500
$ws.Range("L73").Value = 11880
$ws.Range("M73").Value = -3408
$ws.Range("N73").Value = -14064
$ws.Range("H75").Value = 5625
$ws.Range("J75").Value = 6000
$ws.Range("L75").Value = 18000
$ws.Range("N75").Value = -19996
$ws.Range("H76").Value = 5787.857
$ws.Range("I76").Value = 3000
$ws.Range("J76").Value = 6252.5
$ws.Range("K76").Value = 9000
$ws.Range("L76").Value = 18757.5
$ws.Range("M76").Value = -8617
$ws.Range("N76").Value = -19523.5
$ws.Range("H78").Value = 5625
$ws.Range("J78").Value = 6000
$ws.Range("L78").Value = 54000
$ws.Range("N78").Value = -63984
$ws.Range("H79").Value = 5787.857
$ws.Range("I79").Value = 3000
$ws.Range("J79").Value = 6252.5
$ws.Range("K79").Value = 9000
$ws.Range("L79").Value = 18757.5
$ws.Range("M79").Value = -7674
$ws.Range("N79").Value = -21409.5
$ws.Range("H80").Value = 1791.091
$ws.Range("I80").Value = 999.5
$ws.Range("J80").Value = 1967
$ws.Range("K80").Value = 2998.5
$ws.Range("L80").Value = 5901
$ws.Range("M80").Value = -2062.5
$ws.Range("N80").Value = -7773
$ws.Range("H81").Value = 59031.75
$ws.Range("I81").Value = 828.25
$ws.Range("J81").Value = 78432.914
$ws.Range("K81").Value = 2484.75
$ws.Range("L81").Value = 235298.742
$ws.Range("M81").Value = -1361.75
$ws.Range("N81").Value = -237544.742
$ws.Range("H82").Value = 8909.333000000001
$ws.Range("I82").Value = 713
$ws.Range("J82").Value = 13007.5
$ws.Range("K82").Value = 2139
$ws.Range("L82").Value = 39022.5
$ws.Range("M82").Value = -1733
$ws.Range("N82").Value = -39834.5
$ws.Range("H83").Value = 1791.091
$ws.Range("I83").Value = 999.5
$ws.Range("J83").Value = 1967
$ws.Range("K83").Value = 8995.5
$ws.Range("L83").Value = 17703
$ws.Range("M83").Value = -4315.5
$ws.Range("N83").Value = -27063
$ws.Range("H84").Value = 59031.75
$ws.Range("I84").Value = 828.25
$ws.Range("J84").Value = 78432.914
$ws.Range("K84").Value = 7454.25
$ws.Range("L84").Value = 705896.226
$ws.Range("M84").Value = -1838.25
$ws.Range("N84").Value = -717128.226
$ws.Range("H85").Value = 8909.333000000001
$ws.Range("I85").Value = 713
$ws.Range("J85").Value = 13007.5
$ws.Range("K85").Value = 2139
$ws.Range("L85").Value = 39022.5
$ws.Range("M85").Value = -735
$ws.Range("N85").Value = -41830.5
$ws.Range("H135").Value = 563.6
$ws.Range("I135").Value = 559.3570999999999
$ws.Range("J135").Value = 567.3125
$ws.Range("K135").Value = 5034.2139
$ws.Range("L135").Value = 5105.8125
$ws.Range("M135").Value = -2499.2139
$ws.Range("N135").Value = -10175.8125

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").Value = $null
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").Value = $null
